# Fruta / hortaliza, semanal
# Insert 3 new weekly rows (date 2022-01-06 / serial 44578) at the top of the
# Frutilla - Mercado Mayorista Lo Valledor de Santiago dataset, pushing the
# existing rows (771:835) down to (774:838).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 771-835 down by 3 rows (A:T is the full table width)
$ws.Range("A771:T773").Insert()

# New row 771 - Calidad "Especial"
$ws.Range("A771").Value = 6
$ws.Range("B771").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C771").Value = "Metropolitana"
$ws.Range("D771").Value = 44578
$ws.Range("E771").Value = 13
$ws.Range("F771").Value = "Fruta"
$ws.Range("G771").Value = 100101
$ws.Range("H771").Value = "Berries"
$ws.Range("I771").Value = 100112025
$ws.Range("J771").Value = "Frutilla"
$ws.Range("K771").Value = "Sin especificar"
$ws.Range("L771").Value = "Especial"
$ws.Range("M771").Value = 750
$ws.Range("N771").Value = 5000
$ws.Range("O771").Value = 6000
$ws.Range("P771").Value = 5500
$ws.Range("Q771").Value = "`$/bandeja 7 kilos"
$ws.Range("R771").Value = "Provincia de Melipilla"
$ws.Range("S771").Value = 786
$ws.Range("T771").Value = 7

# New row 772 - Calidad "Primera"
$ws.Range("A772").Value = 6
$ws.Range("B772").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C772").Value = "Metropolitana"
$ws.Range("D772").Value = 44578
$ws.Range("E772").Value = 13
$ws.Range("F772").Value = "Fruta"
$ws.Range("G772").Value = 100101
$ws.Range("H772").Value = "Berries"
$ws.Range("I772").Value = 100112025
$ws.Range("J772").Value = "Frutilla"
$ws.Range("K772").Value = "Sin especificar"
$ws.Range("L772").Value = "Primera"
$ws.Range("M772").Value = 600
$ws.Range("N772").Value = 4000
$ws.Range("O772").Value = 5000
$ws.Range("P772").Value = 4500
$ws.Range("Q772").Value = "`$/bandeja 7 kilos"
$ws.Range("R772").Value = "Provincia de Melipilla"
$ws.Range("S772").Value = 643
$ws.Range("T772").Value = 7

# New row 773 - Calidad "Segunda"
$ws.Range("A773").Value = 6
$ws.Range("B773").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C773").Value = "Metropolitana"
$ws.Range("D773").Value = 44578
$ws.Range("E773").Value = 13
$ws.Range("F773").Value = "Fruta"
$ws.Range("G773").Value = 100101
$ws.Range("H773").Value = "Berries"
$ws.Range("I773").Value = 100112025
$ws.Range("J773").Value = "Frutilla"
$ws.Range("K773").Value = "Sin especificar"
$ws.Range("L773").Value = "Segunda"
$ws.Range("M773").Value = 300
$ws.Range("N773").Value = 2500
$ws.Range("O773").Value = 3500
$ws.Range("P773").Value = 3000
$ws.Range("Q773").Value = "`$/bandeja 7 kilos"
$ws.Range("R773").Value = "Provincia de Melipilla"
$ws.Range("S773").Value = 429
$ws.Range("T773").Value = 7
